$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph so the new
# bullet entry can be inserted right after it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $p = $d.Paragraphs.Item($targetIndex)
    $p.Range.InsertParagraphAfter()

    $newP = $d.Paragraphs.Item($targetIndex + 1)
    $newP.Range.Text = "6712818 - Mauricio Lamano Ferreira"
    $newP.Style = "ListBullet"
}
